$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.041.52"
$ws.Range("E2").Value = "  -2.38%  "

$ws.Range("D3").Value = "1.668.03"
$ws.Range("E3").Value = "  -1.81%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "216.83"
$ws.Range("E5").Value = "  -1.42%  "

$ws.Range("D6").Value = "0.5111"
$ws.Range("E6").Value = "  -0.41%  "

$ws.Range("D8").Value = "0.2653"
$ws.Range("E8").Value = "  +0.16%  "

$ws.Range("D9").Value = "0.06405"
$ws.Range("E9").Value = "  +1.61%  "

$ws.Range("D10").Value = "21.87"
$ws.Range("E10").Value = "  -1.55%  "

$ws.Range("D11").Value = "0.07436"
$ws.Range("E11").Value = "  +1.09%  "

$ws.Range("D12").Value = "1.690.35"
$ws.Range("E12").Value = "  -0.53%  "

$ws.Range("D13").Value = "4.498"
$ws.Range("E13").Value = "  -0.48%  "

$ws.Range("D14").Value = "0.5831"
$ws.Range("E14").Value = "  +0.19%  "

$ws.Range("D15").Value = "0.000008547"
$ws.Range("E15").Value = "  +1.22%  "

$ws.Range("D16").Value = "64.32"
$ws.Range("E16").Value = "  -2.10%  "

$ws.Range("D17").Value = "26.095.85"
$ws.Range("E17").Value = "  -2.23%  "

$ws.Range("D18").Value = "4.945"
$ws.Range("E18").Value = "  -1.46%  "

$ws.Range("D20").Value = "10.76"
$ws.Range("E20").Value = "  -2.27%  "

$ws.Range("D21").Value = "190.32"
$ws.Range("E21").Value = "  +1.64%  "

$ws.Range("D22").Value = "6.231"
$ws.Range("E22").Value = "  -0.57%  "

$ws.Range("E23").Value = "  +0.15%  "

$ws.Range("D24").Value = "'145.20"
$ws.Range("E24").Value = "  +0.48%  "

$ws.Range("E25").Value = "  +1.45%  "

$ws.Range("D26").Value = "0.1203"
$ws.Range("E26").Value = "  +3.61%  "

$ws.Range("D27").Value = "15.64"
$ws.Range("E27").Value = "  -0.44%  "

$ws.Range("D28").Value = "0.06576"
$ws.Range("E28").Value = "  +15.98%  "

$ws.Range("D29").Value = "1.327"
$ws.Range("E29").Value = "  -2.03%  "

$ws.Range("D30").Value = "1.317"
$ws.Range("E30").Value = "  -1.37%  "

$ws.Range("D31").Value = "3.545"
$ws.Range("E31").Value = "  +0.99%  "

$ws.Range("D32").Value = "3.518"
$ws.Range("E32").Value = "  +0.82%  "

$ws.Range("D33").Value = "1.645"
$ws.Range("E33").Value = "  +0.21%  "

$ws.Range("E34").Value = "  -0.49%  "

$ws.Range("D35").Value = "0.6098"
$ws.Range("E35").Value = "  +1.07%  "

$ws.Range("D36").Value = "'2.370"
$ws.Range("E36").Value = "  +0.55%  "

$ws.Range("D37").Value = "2.713"
$ws.Range("E37").Value = "  +0.92%  "

$ws.Range("D38").Value = "6.236"
$ws.Range("E38").Value = "  +6.60%  "

$ws.Range("D39").Value = "0.01605"
$ws.Range("E39").Value = "  -0.53%  "

$ws.Range("D40").Value = "1.086.08"
$ws.Range("E40").Value = "  -1.68%  "

$ws.Range("D41").Value = "0.8612"
$ws.Range("E41").Value = "  +0.28%  "

$ws.Range("E42").Value = "  +0.63%  "

$ws.Range("D43").Value = "100.61"
$ws.Range("E43").Value = "  +0.51%  "

$ws.Range("D44").Value = "1.817.43"
$ws.Range("E44").Value = "  -2.16%  "

$ws.Range("E45").Value = "  +4.50%  "

$ws.Range("D46").Value = "56.36"
$ws.Range("E46").Value = "  -0.82%  "

$ws.Range("E47").Value = "  +0.48%  "

$ws.Range("D48").Value = "8.056"
$ws.Range("E48").Value = "  -1.24%  "

$ws.Range("D49").Value = "0.05239"
$ws.Range("E49").Value = "  -0.12%  "

$ws.Range("D51").Value = "5.997"
$ws.Range("E51").Value = "  +3.64%  "
